# model-traits.cfg.xlsx update:
#  - random model formulas (col E) get extra vs(rowf)/vs(columnf)/vs(spl2D(...)) terms
#  - new "label" (col J/K->J) and "label_short" (col K) columns added per-trait
#  - old "mask" (col I) TRUE markers removed
#  - header "description" moves from col J to col L
#  - a couple of cosmetic column-width / selection tweaks

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header row -------------------------------------------------------
$ws.Range("J1").Value = "label"
$ws.Range("K1").Value = "label_short"
$ws.Range("L1").Value = "description"

# --- per-trait label / short-label lookup ------------------------------
$label = @{
    "berry_length"        = "Berry Length"
    "berry_width"         = "Berry Width"
    "berry_weight"        = "Berry Weight"
    "num_seeds"           = "Number of Seeds"
    "num_peds"            = "Number of Pedicels"
    "num_berries"         = "Number of Berries"
    "total_berry_weight"  = "Total Berry Weight"
}
$labelShort = @{
    "berry_length"        = "BL"
    "berry_width"         = "BW"
    "berry_weight"        = "BM"
    "num_seeds"           = "NS"
    "num_peds"            = "NP"
    "num_berries"         = "NB"
    "total_berry_weight"  = "TBM"
}

# --- data rows 2..22 : per-year models (random -> vs(id,Gu=A)+...) -----
for ($r = 2; $r -le 22; $r++) {
    $ws.Cells.Item($r, 5).Value = "~vs(id, Gu=A)+vs(rowf)+vs(columnf)+vs(spl2D(row,column))"
    $trait = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 10).Value = $label[$trait]
    $ws.Cells.Item($r, 11).Value = $labelShort[$trait]
}

# --- mask column (I) TRUE markers removed on rows 6,7,14,17,21,22 ------
$ws.Range("I6").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("I14").ClearContents()
$ws.Range("I17").ClearContents()
$ws.Range("I21").ClearContents()
$ws.Range("I22").ClearContents()

# --- data rows 23..29 : all-years models (random -> ... + id:year + ...) --
for ($r = 23; $r -le 29; $r++) {
    $ws.Cells.Item($r, 5).Value = "~vs(id, Gu=A) + id:year + vs(rowf)+vs(columnf)+vs(spl2D(row,column))"
    $trait = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 10).Value = $label[$trait]
    $ws.Cells.Item($r, 11).Value = $labelShort[$trait]
}

# --- cosmetic: column widths for the now-wider / new columns -----------
$ws.Columns.Item(5).ColumnWidth = 60.330729166666664
$ws.Columns.Item(7).ColumnWidth = 9.830729166666666
$ws.Columns.Item(9).ColumnWidth = 4.666666666666667
$ws.Columns.Item(10).ColumnWidth = 16.498697916666668
$ws.Columns.Item(11).ColumnWidth = 16.498697916666668
$ws.Columns.Item(12).ColumnWidth = 62.498697916666664

# --- selection moves to K29 --------------------------------------------
$ws.Range("K29").Select()
